# "catalogo completo - parcial"
# Rearranges the "Right Arrow" step chevrons on slide 3 of the
# "SKU to VTEXADMIN to SITE" deck: a few existing steps move up/down to
# make room for two brand-new steps ("CAMPO DE PRODUTO" and
# "ESPECIFICAÇÃO DE SKU"), the old "ESPECIFICAÇÃO DE PRODUTO" / "IMAGEM DE
# SKU" steps are re-homed further down the chain, the "ATIVA PRODUTO" step
# is dropped altogether, and the step that used to read "ESPECIFICAÇÃO DE
# SKU" is renamed to "CAMPO DE SKU".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        if ($shapes.Item($i).Id -eq $id) {
            return $shapes.Item($i)
        }
    }
    return $null
}

# PowerPoint's COM surface works in points (Shape.Left/.Top/.Width/.Height)
# while the OOXML on disk stores EMU (1 pt = 12700 EMU), and this host
# rounds Left/Top through a single-precision float on the way down, so a
# naive "emu / 12700" literal can truncate one EMU short once it comes
# back out. Nudge the point value up in tiny steps until the round trip
# lands exactly back on the EMU value we actually want.
function Emu2Pt($targetEmu) {
    $base = [double]$targetEmu / 12700.0
    for ($i = 0; $i -lt 2000; $i++) {
        $candidate = $base + ($i * 0.0000001)
        $f32 = [float]$candidate
        $backEmu = [int64]([double]$f32 * 12700.0)
        if ($backEmu -eq $targetEmu) {
            return $candidate
        }
    }
    return $base
}

# ---------------------------------------------------------------------
# 1. Nudge the steps that stay put, up/down a little, to free up space.
# ---------------------------------------------------------------------
(Get-ShapeById $s.Shapes 85).Top = Emu2Pt 516383   # CATEGORIA
(Get-ShapeById $s.Shapes 87).Top = Emu2Pt 837515   # MARCA
(Get-ShapeById $s.Shapes 91).Top = Emu2Pt 1479765  # PRODUTO
(Get-ShapeById $s.Shapes 93).Top = Emu2Pt 2609156  # SKU
(Get-ShapeById $s.Shapes 96).Top = Emu2Pt 4217519  # PREÇO
(Get-ShapeById $s.Shapes 97).Top = Emu2Pt 4663837  # ESTOQUE
(Get-ShapeById $s.Shapes 102).Top = Emu2Pt 5519058 # ATIVA SKU

# ---------------------------------------------------------------------
# 2. Rename the "ESPECIFICAÇÃO DE SKU" step to "CAMPO DE SKU" (still two
#    runs, same split point as before: "CAMPO DE " + "SKU").
# ---------------------------------------------------------------------
$shp94 = Get-ShapeById $s.Shapes 94
$tr94 = $shp94.TextFrame.TextRange
$tr94.Text = "CAMPO DE "
$tr94.InsertAfter("SKU") | Out-Null

# ---------------------------------------------------------------------
# 3. Drop the "ATIVA PRODUTO" step entirely (not replaced by anything).
# ---------------------------------------------------------------------
(Get-ShapeById $s.Shapes 100).Delete()

# ---------------------------------------------------------------------
# 4. Add the two brand new steps, "CAMPO DE PRODUTO" and "ESPECIFICAÇÃO
#    DE SKU", by duplicating a similar chevron and re-texting it.
# ---------------------------------------------------------------------
$tmplA = Get-ShapeById $s.Shapes 93
$dup = $tmplA.Duplicate()
$new119 = $dup.Item(1)
$new119.Left = Emu2Pt 3025518
$new119.Top = Emu2Pt 1754625
$new119.TextFrame.TextRange.Text = "CAMPO DE PRODUTO"

$tmplB = Get-ShapeById $s.Shapes 96
$dup = $tmplB.Duplicate()
$new120 = $dup.Item(1)
$new120.Left = Emu2Pt 3033682
$new120.Top = Emu2Pt 3207873
$new120.TextFrame.TextRange.Text = "ESPECIFICAÇÃO DE SKU"

# ---------------------------------------------------------------------
# 5. Duplicate the two steps that get moved further down the chain
#    ("IMAGEM DE SKU" and "ESPECIFICAÇÃO DE PRODUTO") before deleting
#    their old copies, so the duplicates inherit the exact same style.
# ---------------------------------------------------------------------
$old95 = Get-ShapeById $s.Shapes 95
$dup = $old95.Duplicate()
$new95 = $dup.Item(1)
$new95.Left = Emu2Pt 3020070
$new95.Top = Emu2Pt 3518118
$old95.Delete()

$old92 = Get-ShapeById $s.Shapes 92
$dup = $old92.Duplicate()
$new92 = $dup.Item(1)
$new92.Left = Emu2Pt 3020070
$new92.Top = Emu2Pt 2048539
$old92.Delete()
